$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source added a new weekly price sample. In the source data, rows are
# ordered newest-first, so the new observation is inserted at row 161,
# pushing the former rows 161-168 down to 162-169 (dimension grows to R169).
$ws.Rows.Item(161).Insert()

$newRow = 161
$ws.Cells.Item($newRow, 1).Value  = 8
$ws.Cells.Item($newRow, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item($newRow, 3).Value  = "Coquimbo"
$ws.Cells.Item($newRow, 4).Value  = 45147
$ws.Cells.Item($newRow, 5).Value  = 4
$ws.Cells.Item($newRow, 6).Value  = 100114007
$ws.Cells.Item($newRow, 7).Value  = "Jengibre"
$ws.Cells.Item($newRow, 8).Value  = "Sin especificar"
$ws.Cells.Item($newRow, 9).Value  = "Primera"
$ws.Cells.Item($newRow, 10).Value = 400
$ws.Cells.Item($newRow, 11).Value = 18000
$ws.Cells.Item($newRow, 12).Value = 19000
$ws.Cells.Item($newRow, 13).Value = 18500
$ws.Cells.Item($newRow, 14).Value = "`$/caja 13 kilos"
$ws.Cells.Item($newRow, 15).Value = "Perú"
$ws.Cells.Item($newRow, 16).Value = 1423
$ws.Cells.Item($newRow, 17).Value = 13
$ws.Cells.Item($newRow, 18).Value = "Hortaliza"
